# [Group18] Danh Gia Thanh Vien Lan 6
# Fill in "Round 6 evaluation" scores (column K) for 4 members (rows 16-19)
# in the "Ghi cong" sheet, then leave the selection on K20 as the last
# active cell (mirrors what a user would do after typing values down
# column K and landing on the next empty row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Ghi cong")

$ws.Range("K16").Value = 1
$ws.Range("K17").Value = 3
$ws.Range("K18").Value = 1
$ws.Range("K19").Value = 3

$ws.Activate()
$ws.Range("K20").Select()
